# Auto-generated edit script applying the cell-value changes described in the commit diff.
# Each worksheet is selected by name, then cells are updated to their new values.
# Cells that were removed entirely in the diff are cleared with ClearContents().

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4855.026
$ws.Range("I15").Value = 4855.026
$ws.Range("K15").Value = 14565.078
$ws.Range("M15").Value = -14396.078
$ws.Range("H17").Value = 2843.7693
$ws.Range("J17").Value = 2843.7693
$ws.Range("L17").Value = 8531.3079
$ws.Range("N17").Value = -8867.3079
$ws.Range("H32").Value = 7737.25
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652
$ws.Range("H33").Value = 17369.842
$ws.Range("J33").Value = 956
$ws.Range("L33").Value = 956
$ws.Range("N33").Value = -1414
$ws.Range("H51").Value = 7300
$ws.Range("I51").Value = 8275
$ws.Range("K51").Value = 8275
$ws.Range("M51").Value = -7791
$ws.Range("H62").Value = 4166.6665
$ws.Range("J62").Value = 3999
$ws.Range("L62").Value = 3999
$ws.Range("N62").Value = -5247
$ws.Range("H65").Value = 4166.6665
$ws.Range("J65").Value = 3999
$ws.Range("L65").Value = 19995
$ws.Range("N65").Value = -26235
$ws.Range("H96").Value = 354.33334
$ws.Range("I96").Value = 388.5
$ws.Range("J96").Value = 217.66667
$ws.Range("K96").Value = 1165.5
$ws.Range("L96").Value = 653.00001
$ws.Range("M96").Value = 207.5
$ws.Range("N96").Value = -3399.00001
$ws.Range("H103").Value = 2570.2
$ws.Range("I103").Value = 5226.5
$ws.Range("J103").Value = 799.3333
$ws.Range("K103").Value = 15679.5
$ws.Range("L103").Value = 2397.9999
$ws.Range("M103").Value = -15093.5
$ws.Range("N103").Value = -3569.9999
$ws.Range("H132").Value = 1613.3243
$ws.Range("I132").Value = 1528.5625
$ws.Range("K132").Value = 4585.6875
$ws.Range("M132").Value = -2055.6875
$ws.Range("H137").Value = 1065.0312
$ws.Range("I137").Value = 996.1070999999999
$ws.Range("K137").Value = 2988.3213
$ws.Range("M137").Value = -438.3212999999996
$ws.Range("H138").Value = 3504.7407
$ws.Range("I138").Value = 2083.2856
$ws.Range("J138").Value = 4002.25
$ws.Range("K138").Value = 6249.8568
$ws.Range("L138").Value = 12006.75
$ws.Range("M138").Value = -1109.8568
$ws.Range("N138").Value = -22286.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4881.4224
$ws.Range("J61").Value = 4998
$ws.Range("L61").Value = 4998
$ws.Range("N61").Value = -5422
$ws.Range("H74").Value = 3480.175
$ws.Range("I74").Value = 1769.909
$ws.Range("K74").Value = 1769.909
$ws.Range("M74").Value = -895.9090000000001
$ws.Range("H77").Value = 3480.175
$ws.Range("I77").Value = 1769.909
$ws.Range("K77").Value = 8849.545
$ws.Range("M77").Value = -4481.545
$ws.Range("H122").Value = 1936.909
$ws.Range("I122").Value = 1936.909
$ws.Range("K122").Value = 5810.727000000001
$ws.Range("M122").Value = -3360.727000000001
$ws.Range("H132").Value = 1556.317
$ws.Range("I132").Value = 1556.317
$ws.Range("K132").Value = 4668.951
$ws.Range("M132").Value = -2138.951
$ws.Range("H136").Value = 4881.4224
$ws.Range("J136").Value = 4998
$ws.Range("L136").Value = 14994
$ws.Range("N136").Value = -20094

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4238.6665
$ws.Range("I99").Value = 3066.647
$ws.Range("K99").Value = 3066.647
$ws.Range("M99").Value = -1568.647

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 19954.5
$ws.Range("J69").Value = 29910
$ws.Range("L69").Value = 29910
$ws.Range("N69").Value = -31408
$ws.Range("H72").Value = 19954.5
$ws.Range("J72").Value = 29910
$ws.Range("L72").Value = 89730
$ws.Range("N72").Value = -97218
$ws.Range("H132").Value = 5024.9165
$ws.Range("J132").Value = 6666.5
$ws.Range("L132").Value = 19999.5
$ws.Range("N132").Value = -25059.5
$ws.Range("H134").Value = 4335.263
$ws.Range("I134").Value = 3080.8333
$ws.Range("J134").Value = 6485.7144
$ws.Range("K134").Value = 9242.499899999999
$ws.Range("L134").Value = 19457.1432
$ws.Range("M134").Value = -6707.499899999999
$ws.Range("N134").Value = -24527.1432

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 78.40000000000001
$ws.Range("I60").Value = 78
$ws.Range("J60").Value = 80
$ws.Range("K60").Value = 234
$ws.Range("L60").Value = 240
$ws.Range("M60").Value = 17
$ws.Range("N60").Value = -742
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242
$ws.Range("H107").Value = 449.75
$ws.Range("J107").Value = 449.75
$ws.Range("L107").Value = 1349.25
$ws.Range("N107").Value = -5189.25
$ws.Range("H124").Value = 2499
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H139").Value = 7540.913
$ws.Range("I139").Value = 3503.5625
$ws.Range("K139").Value = 10510.6875
$ws.Range("M139").Value = -5370.6875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3720.889
$ws.Range("I40").Value = 3586.125
$ws.Range("J40").Value = 4799
$ws.Range("K40").Value = 3586.125
$ws.Range("L40").Value = 4799
$ws.Range("M40").Value = -3450.125
$ws.Range("N40").Value = -5071
$ws.Range("H46").Value = 3682.4375
$ws.Range("I46").Value = 2224.5386
$ws.Range("K46").Value = 2224.5386
$ws.Range("M46").Value = -2036.5386
$ws.Range("H61").Value = 28560.838
$ws.Range("I61").Value = 41810.08
$ws.Range("J61").Value = 958.25
$ws.Range("K61").Value = 41810.08
$ws.Range("L61").Value = 958.25
$ws.Range("M61").Value = -41608.08
$ws.Range("N61").Value = -1362.25
$ws.Range("H100").Value = 3891.7727
$ws.Range("I100").Value = 2147.182
$ws.Range("K100").Value = 2147.182
$ws.Range("M100").Value = -1606.182
$ws.Range("H113").Value = 28560.838
$ws.Range("I113").Value = 41810.08
$ws.Range("J113").Value = 958.25
$ws.Range("K113").Value = 41810.08
$ws.Range("L113").Value = 958.25
$ws.Range("M113").Value = -39640.08
$ws.Range("N113").Value = -5298.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 6996
$ws.Range("H100").Value = 826.38464
$ws.Range("I100").Value = 704.7778
$ws.Range("K100").Value = 1409.5556
$ws.Range("M100").Value = -868.5555999999999
$ws.Range("H107").Value = 1793.7059
$ws.Range("J107").Value = 2617.5715
$ws.Range("L107").Value = 7852.7145
$ws.Range("N107").Value = -11692.7145
$ws.Range("H122").Value = 3089.5217
$ws.Range("I122").Value = 1861.3572
$ws.Range("K122").Value = 5584.071599999999
$ws.Range("M122").Value = -3134.071599999999
$ws.Range("H126").Value = 6659.222
$ws.Range("I126").Value = 4181.615
$ws.Range("K126").Value = 12544.845
$ws.Range("M126").Value = -10074.845
$ws.Range("H136").Value = 5260.522
$ws.Range("I136").Value = 4841.7896
$ws.Range("K136").Value = 14525.3688
$ws.Range("M136").Value = -11975.3688
